$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A67").Value = "Andrea Conzatti"
$ws.Range("B67").Value = "Daniele Dalbosco | IMONTAGNA"
$ws.Range("C67").Value = "Luca Frasca | Clitoriders"
$ws.Range("D67").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("E67").Value = "Geremia Carollo | FC SAVIGNANO"
$ws.Range("F67").Value = "Davide Simoncelli | Avanzi"
